$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45184

for ($row = 2; $row -le 498; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = $newDate
    }
}
